# Word COM-interop script implementing the commit "unnatrual things with gams"
# against scaffold/Writing/abstract.docx.
#
# Strategy: perform a sequence of surgical Find & Replace operations on
# $d.Content (and on narrowed Range objects where an operation must target a
# specific occurrence of otherwise-duplicated text), keeping every
# find/replace pair inside a single run's formatting span wherever the
# replacement text must preserve (or must avoid acquiring) italics.

$d = $word.ActiveDocument
$cr = [char]13
$wdFindContinue = 1
$wdReplaceOne = 1
$wdReplaceAll = 2

function Replace-One($text, $replacement) {
    $ok = $d.Content.Find.Execute($text, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replacement, $wdReplaceOne)
    if (-not $ok) {
        throw "Find/Replace failed (ReplaceOne) for: $text"
    }
}

function Replace-All($text, $replacement) {
    $ok = $d.Content.Find.Execute($text, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replacement, $wdReplaceAll)
    if (-not $ok) {
        throw "Find/Replace failed (ReplaceAll) for: $text"
    }
}

# ---------------------------------------------------------------------------
# 1) "Background and Methods" paragraph: the near-complete compensation
#    timing phrase changes from "in the 1990s following" to
#    "beginning in 1996, with".
# ---------------------------------------------------------------------------
Replace-All `
    ". Partial energetic compensation from smaller granivores was observed immediately, and near-complete compensation occurred in the 1990s following the establishment of the " `
    ". Partial energetic compensation from smaller granivores was observed immediately, and near-complete compensation occurred beginning in 1996, with the establishment of the "

# ---------------------------------------------------------------------------
# 2) "Results and Conclusions" paragraph rewrite.
#    Split into pieces that never cross the italic "C. baileyi" runs, so the
#    italic formatting on the surviving "C. baileyi" mention is preserved
#    untouched, and the text around the second (deleted) mention stays
#    plain/non-italic.
# ---------------------------------------------------------------------------

# 2a. Opening clause, up to (not including) the first "C. baileyi".
Replace-One `
    "Since 2010, compensatory gains in energy use from small granivores on treatment plots relative to controls declined to near zero, contrasting to the partial compensation prior to the establishment of " `
    "Since 2010, total energy use on kangaroo rat exclosure plots declined to ~40% that on controls, compared to 70% from 1996-2010 and 24% before 1996. This coincided with a precipitous decline in "

# 2b. Middle clause, strictly between the first and second "C. baileyi"
#     mentions (does not touch either italic run).
Replace-One `
    " and the near-complete compensation following its arrival. This coincided with a long-term, sitewide increase in the proportion of energy use from smaller granivores, and a recent precipitous decline of " `
    ", a long-term increase in the proportion of energy use from small granivores sitewide – from 5% of total energy use on control plots before 1996, to 30% after 2010 – and decreasing gains in small granivore energy use on treatment relative to control plots – from a threefold increase before 1996, to near-convergence after 2010"

# 2c. Delete the now-unneeded second "C. baileyi" italic mention outright
#     (match is wholly inside the italic run, so surrounding plain-text runs
#     join cleanly with no stray formatting or extra spaces).
Replace-One `
    "C. baileyi. Therefore, while kangaroo rat removal" `
    ". Therefore, while kangaroo rat removal"

# Wait: the line above would start the match at the italic run and bleed
# italics onto the replacement. Instead remove just the italic words.
